$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder countries (Bosnia y Herzegovina now listed before Principado de Andorra) ---
$ws.Range("A79").Value = "Bosnia y Herzegovina"
$ws.Range("A80").Value = "Principado de Andorra"

# --- Reorder countries (Honduras now listed right after Puerto Rico) ---
$ws.Range("A118").Value = "Honduras"
$ws.Range("A119").Value = "Guam"
$ws.Range("A120").Value = "Bolivia"

# --- Update "last updated" timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 01:16"

# --- Estados Unidos (row 6) ---
$ws.Range("B6").Value = 43537
$ws.Range("C6").Value = 9971
$ws.Range("E6").Value = 42697

# --- Uruguay (row 75) ---
$ws.Range("B75").Value = 162
$ws.Range("C75").Value = 4
$ws.Range("E75").Value = 162
$ws.Range("F75").Value = 3

# --- Bosnia y Herzegovina, now row 79 ---
$ws.Range("B79").Value = 136
$ws.Range("C79").Value = 10
$ws.Range("D79").Value = 2
$ws.Range("E79").Value = 133
$ws.Range("F79").Value = 1

# --- Principado de Andorra, now row 80 ---
$ws.Range("B80").Value = 133
$ws.Range("C80").Value = 20
$ws.Range("D80").Value = 1
$ws.Range("E80").Value = 131
$ws.Range("F80").Value = 2

# --- Honduras, now row 118 ---
$ws.Range("B118").Value = 30
$ws.Range("C118").Value = 4
$ws.Range("E118").Value = 30
$ws.Range("H118").Value = 0

# --- Guam, now row 119 ---
$ws.Range("B119").Value = 29
$ws.Range("C119").Value = 2
$ws.Range("E119").Value = 28
$ws.Range("H119").Value = 1

# --- Bolivia, now row 120 ---
$ws.Range("C120").Value = 3
